# Update the "取得日時" (acquired datetime) timestamp in column A
# for all data rows (2-20) on the first worksheet from
# "2025-10-01 18:24:11" to "2025-10-01 18:32:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-01 18:32:00"

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
